# Append two new data rows (173 and 174) to the VLC.MI price-history sheet,
# mirroring results produced by the R script that feeds this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 173 ----
$rowA = $ws.Cells.Item(173, 1)
$ws.Cells.Item(172, 1).Copy($rowA)      # reuse the date/time cell style (s="1")
$rowA.Value = 45447.2916666667
$ws.Cells.Item(173, 2).Value = 0
$ws.Cells.Item(173, 3).Value = 5.30000019073486
$ws.Cells.Item(173, 4).Value = 5.30000019073486
$ws.Cells.Item(173, 5).Value = 5.30000019073486
$ws.Cells.Item(173, 6).Value = 5.30000019073486

$rowG = $ws.Cells.Item(173, 7)
$rowG.NumberFormat = "@"                # force text storage so it lands in sharedStrings
$rowG.Value = "5.30000019073486"
$rowG.Style = "Normal"                  # drop the temporary text style again

$ws.Cells.Item(173, 8).Value = "VLC.MI"

# ---- Row 174 ----
$rowA2 = $ws.Cells.Item(174, 1)
$ws.Cells.Item(172, 1).Copy($rowA2)
$rowA2.Value = 45448.5777546296
$ws.Cells.Item(174, 2).Value = 1620
$ws.Cells.Item(174, 3).Value = 5.15000009536743
$ws.Cells.Item(174, 4).Value = 4.90000009536743
$ws.Cells.Item(174, 5).Value = 5.09999990463257
$ws.Cells.Item(174, 6).Value = 5.15000009536743

$rowG2 = $ws.Cells.Item(174, 7)
$rowG2.NumberFormat = "@"
$rowG2.Value = "5.15000009536743"
$rowG2.Style = "Normal"

$ws.Cells.Item(174, 8).Value = "VLC.MI"
